$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2 = @{ "G" = 2.629732666666667; "H" = 7.889198; "I" = 0.07156737804735891; "J" = 0.07156737804735891; "M" = 133.7780026666667; "N" = 401.334008; "O" = 0.50863533211804; "P" = 0.5086353321180399; "Q" = 351.8003836939538; "R" = 3166.203453245585; "S" = 0.03640169710193573; "T" = 0.03640169710193572 }
    3 = @{ "G" = 2.629732666666667; "H" = 7.889198; "I" = 0.07156737804735891; "J" = 0.07156737804735891; "O" = 0.1993888292903622; "P" = 0.1993888292903622; "Q" = 137.908364243086; "R" = 1241.175278187774; "S" = 0.01426973572424366; "T" = 0.01426973572424366 }
    4 = @{ "G" = 2.629732666666667; "H" = 7.889198; "I" = 0.07156737804735891; "J" = 0.07156737804735891; "M" = 21.197691; "N" = 63.593073; "O" = 0.08059542216956049; "P" = 0.08059542216956046; "Q" = 55.74426048060601; "R" = 501.6983443254541; "S" = 0.005768003047295427; "T" = 0.005768003047295425 }
    5 = @{ "G" = 2.629732666666667; "H" = 7.889198; "I" = 0.07156737804735891; "J" = 0.07156737804735891; "M" = 55.59592133333333; "N" = 166.787764; "O" = 0.2113804164220374; "P" = 0.2113804164220373; "Q" = 146.2024104636969; "R" = 1315.821694173272; "S" = 0.0151279421738841; "T" = 0.0151279421738841 }
    6 = @{ "I" = 0.493312042610523; "J" = 0.493312042610523; "M" = 133.7780026666667; "N" = 401.334008; "O" = 0.50863533211804; "P" = 0.5086353321180399; "Q" = 2424.950733229141; "R" = 21824.55659906227; "S" = 0.2509159346310321; "T" = 0.250915934631032 }
    7 = @{ "I" = 0.493312042610523; "J" = 0.493312042610523; "O" = 0.1993888292903622; "P" = 0.1993888292903622; "Q" = 950.5987045216812; "R" = 8555.38834069513; "S" = 0.09836091065094947; "T" = 0.09836091065094946 }
    8 = @{ "I" = 0.493312042610523; "J" = 0.493312042610523; "M" = 21.197691; "N" = 63.593073; "O" = 0.08059542216956049; "P" = 0.08059542216956046; "Q" = 384.2437120346011; "R" = 3458.19340831141; "S" = 0.03975869233552332; "T" = 0.0397586923355233 }
    9 = @{ "I" = 0.493312042610523; "J" = 0.493312042610523; "M" = 55.59592133333333; "N" = 166.787764; "O" = 0.2113804164220374; "P" = 0.2113804164220373; "Q" = 1007.76934559069; "R" = 9069.924110316211; "S" = 0.1042765049930182; "T" = 0.1042765049930182 }
    10 = @{ "G" = 7.550656333333333; "H" = 22.651969; "I" = 0.2054888252189962; "J" = 0.2054888252189962; "M" = 133.7780026666667; "N" = 401.334008; "O" = 0.50863533211804; "P" = 0.5086353321180399; "Q" = 1010.11172309575; "R" = 9091.005507861753; "S" = 0.10451887686181; "T" = 0.10451887686181 }
    11 = @{ "G" = 7.550656333333333; "H" = 22.651969; "I" = 0.2054888252189962; "J" = 0.2054888252189962; "O" = 0.1993888292903622; "P" = 0.1993888292903622; "Q" = 395.971300463633; "R" = 3563.741704172698; "S" = 0.04097217629266751; "T" = 0.04097217629266751 }
    12 = @{ "G" = 7.550656333333333; "H" = 22.651969; "I" = 0.2054888252189962; "J" = 0.2054888252189962; "M" = 21.197691; "N" = 63.593073; "O" = 0.08059542216956049; "P" = 0.08059542216956046; "Q" = 160.056479801193; "R" = 1440.508318210737; "S" = 0.01656145861965203; "T" = 0.01656145861965202 }
    13 = @{ "G" = 7.550656333333333; "H" = 22.651969; "I" = 0.2054888252189962; "J" = 0.2054888252189962; "M" = 55.59592133333333; "N" = 166.787764; "O" = 0.2113804164220374; "P" = 0.2113804164220373; "Q" = 419.7856955230351; "R" = 3778.071259707316; "S" = 0.04343631344486667; "T" = 0.04343631344486666 }
    14 = @{ "G" = 8.437784666666667; "H" = 25.313354; "I" = 0.2296317541231219; "J" = 0.2296317541231219; "M" = 133.7780026666667; "N" = 401.334008; "O" = 0.50863533211804; "P" = 0.5086353321180399; "Q" = 1128.789979638093; "R" = 10159.10981674283; "S" = 0.1167988235232622; "T" = 0.1167988235232622 }
    15 = @{ "G" = 8.437784666666667; "H" = 25.313354; "I" = 0.2296317541231219; "J" = 0.2296317541231219; "O" = 0.1993888292903622; "P" = 0.1993888292903622; "Q" = 442.4940587935781; "R" = 3982.446529142203; "S" = 0.04578600662250158; "T" = 0.04578600662250157 }
    16 = @{ "G" = 8.437784666666667; "H" = 25.313354; "I" = 0.2296317541231219; "J" = 0.2296317541231219; "M" = 21.197691; "N" = 63.593073; "O" = 0.08059542216956049; "P" = 0.08059542216956046; "Q" = 178.861552088538; "R" = 1609.753968796842; "S" = 0.01850726816708972; "T" = 0.01850726816708971 }
    17 = @{ "G" = 8.437784666666667; "H" = 25.313354; "I" = 0.2296317541231219; "J" = 0.2296317541231219; "M" = 55.59592133333333; "N" = 166.787764; "O" = 0.2113804164220374; "P" = 0.2113804164220373; "Q" = 469.1064125556063; "R" = 4221.957713000455; "S" = 0.0485396558102684; "T" = 0.04853965581026838 }
}

foreach ($row in $changes.Keys) {
    $rowData = $changes[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
